$d = $word.ActiveDocument

# --- Change 1: add " (draft file is in icat3-reporting > docs)" after "Images.properties" ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Images.properties", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $ins1 = $d.Range($rng1.End, $rng1.End)
    $ins1.InsertAfter(" (draft file is in icat3-reporting > docs)")
}

# --- Change 2: add three new list paragraphs (dbUsername, dbPassword, dbURL) after the
# "sourceFolder" bullet item (before the blank line that precedes "Libraries:") ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("sourceFolder", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $srcPara = $rng2.Paragraphs(1)
    $srcRange = $srcPara.Range
    $srcRange.InsertParagraphAfter()
    $dbUserPara = $srcPara.Next()
    $r = $dbUserPara.Range
    $r.Collapse(1)
    $r.InsertAfter("dbUsername")
    $r.Collapse(0)
    $r.InsertAfter("- username for connecting to the logging database")

    $dbUserRange = $dbUserPara.Range
    $dbUserRange.InsertParagraphAfter()
    $dbPassPara = $dbUserPara.Next()
    $r2 = $dbPassPara.Range
    $r2.Collapse(1)
    $r2.InsertAfter("dbPassword")
    $r2.Collapse(0)
    $r2.InsertAfter("- password for connecting to the logging database")

    $dbPassRange = $dbPassPara.Range
    $dbPassRange.InsertParagraphAfter()
    $dbUrlPara = $dbPassPara.Next()
    $r3 = $dbUrlPara.Range
    $r3.Collapse(1)
    $r3.InsertAfter("dbURL")
    $r3.Collapse(0)
    $r3.InsertAfter("- url for connecting to the logging database")
}
